# The deck shipped with two theme parts:
#   ppt/theme/theme1.xml -> "Integral"     (wired to the slide master / used by all slides)
#   ppt/theme/theme2.xml -> "Office Theme" (wired to the notes master only)
#
# The authored change swaps the two themes' content: the slide master (and
# therefore every slide) switches from the "Integral" palette to the
# standard "Office" palette, while the old "Integral" palette moves over to
# become the (otherwise unused) notes-master theme.
#
# The only real differences between the two theme parts are the <a:theme>
# name, the <a:clrScheme> name and its twelve colours (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) -- font scheme and format scheme are byte for
# byte identical between the two themes, so nothing else needs to change.

function ConvertTo-OleColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Target palette: the "Office Theme" colours that used to live in theme2.xml.
# Order matches the ThemeColorScheme / ColorScheme index order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
#   9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

# Apply the new palette through every slide's theme-colour scheme (there is
# a single slide master/theme in this deck, so any slide reaches it).
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = ConvertTo-OleColor($officeColors[$i - 1])
}
